$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct the Week 4 dates (they were duplicated from Week 3) ---
$ws.Range("A52").Value = 42870
$ws.Range("A52").NumberFormat = "m/d/yy"

$ws.Range("A55").Value = 42871
$ws.Range("A55").NumberFormat = "m/d/yy"

$ws.Range("A58").Value = 42872
$ws.Range("A58").NumberFormat = "m/d/yy"

$ws.Range("A61").Value = 42873
$ws.Range("A61").NumberFormat = "m/d/yy"
$ws.Range("I61").Value = "ziek"

$ws.Range("A64").Value = 42874
$ws.Range("A64").NumberFormat = "m/d/yy"
$ws.Range("I64").Value = "ziek"

# --- Style (yellow) separator row under the Week 4 block ---
$ws.Range("A66:K66").Interior.Color = 65535

# --- Week 5 block ---
# Monday
$ws.Range("A67").Value = 42884
$ws.Range("A67").NumberFormat = "m/d/yy"
$ws.Range("C67").Value = "ma"
$ws.Range("E67").NumberFormat = "h:mm"
$ws.Range("G67").Value = "School"
$ws.Range("I67").Value = "C# app gewerkt"

$ws.Range("G68").Value = "School"

# Tuesday
$ws.Range("A70").Value = 42885
$ws.Range("A70").NumberFormat = "m/d/yy"
$ws.Range("C70").Value = "di"
$ws.Range("E70").NumberFormat = "h:mm"
$ws.Range("G70").Value = "School"
$ws.Range("I70").Value = "controleren documentatie"

$ws.Range("G71").Value = "School"

# Now label the whole block as "Week 5" (first appearance of this string)
$ws.Range("B67").Value = "Week 5"
$ws.Range("B68").Value = "Week 5"
$ws.Range("B70").Value = "Week 5"
$ws.Range("B71").Value = "Week 5"

# Wednesday
$ws.Range("A73").Value = 42886
$ws.Range("A73").NumberFormat = "m/d/yy"
$ws.Range("B73").Value = "Week 5"
$ws.Range("C73").Value = "wo"
$ws.Range("E73").NumberFormat = "h:mm"
$ws.Range("G73").Value = "School"

$ws.Range("B74").Value = "Week 5"
$ws.Range("G74").Value = "School"

# Thursday
$ws.Range("A76").Value = 42887
$ws.Range("A76").NumberFormat = "m/d/yy"
$ws.Range("B76").Value = "Week 5"
$ws.Range("C76").Value = "do"
$ws.Range("E76").NumberFormat = "h:mm"
$ws.Range("G76").Value = "School"

$ws.Range("B77").Value = "Week 5"
$ws.Range("G77").Value = "School"

# Friday
$ws.Range("A79").Value = 42888
$ws.Range("A79").NumberFormat = "m/d/yy"
$ws.Range("B79").Value = "Week 5"
$ws.Range("C79").Value = "vr"
$ws.Range("E79").NumberFormat = "h:mm"
$ws.Range("G79").Value = "School"

$ws.Range("B80").Value = "Week 5"
$ws.Range("G80").Value = "School"

# --- Style (yellow) separator row closing the Week 5 block ---
$ws.Range("A81:K81").Interior.Color = 65535

# --- Restore the view: scroll position + selection on the new last entry ---
[void]$ws.Range("C79").Select()
